$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 5159
$ws.Range("F5").Value = 7458
$ws.Range("F11").Value = 29
$ws.Range("F12").Value = 4322
$ws.Range("F13").Value = 1762
$ws.Range("F14").Value = 104
$ws.Range("F15").Value = 107
$ws.Range("F17").Value = 582
$ws.Range("F19").Value = 207
$ws.Range("F20").Value = 502
$ws.Range("F21").Value = 440
$ws.Range("F23").Value = 308
$ws.Range("F24").Value = 101
$ws.Range("F26").Value = 1187
$ws.Range("F28").Value = 1379
$ws.Range("F29").Value = 107
$ws.Range("F31").Value = 28
$ws.Range("F34").Value = 62
$ws.Range("F36").Value = 66
$ws.Range("F37").Value = 2902
$ws.Range("F38").Value = 707
$ws.Range("F39").Value = 21
$ws.Range("F40").Value = 78
$ws.Range("F42").Value = 31
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 12
$ws.Range("F3").Value = 11
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5159
$ws.Range("F5").Value = 7458
$ws.Range("F11").Value = 29
$ws.Range("F12").Value = 4322
$ws.Range("F13").Value = 1762
$ws.Range("F14").Value = 104
$ws.Range("F15").Value = 107
$ws.Range("F17").Value = 582
$ws.Range("F19").Value = 207
$ws.Range("F20").Value = 502
$ws.Range("F21").Value = 440
$ws.Range("F23").Value = 12
$ws.Range("F24").Value = 308
$ws.Range("F25").Value = 101
$ws.Range("F27").Value = 1187
$ws.Range("F29").Value = 1379
$ws.Range("F30").Value = 107
$ws.Range("F32").Value = 28
$ws.Range("F35").Value = 62
$ws.Range("F37").Value = 66
$ws.Range("F38").Value = 2902
$ws.Range("F39").Value = 11
$ws.Range("F40").Value = 707
$ws.Range("F41").Value = 21
$ws.Range("F42").Value = 78
$ws.Range("F44").Value = 31
